# fix bug exeded requeste in google drive
# Refresh the price list date stamp and the unit prices that were pulled
# from the (rate-limited) Google Drive price source.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Bump the date stamp in A1 by one day (serial 45310 -> 45311)
$ws.Range("A1").Value = 45311

# Updated unit prices in column D (rows 25-30), recalculated from the
# refreshed source data
$ws.Range("D25").Value = 635.976
$ws.Range("D26").Value = 709.125
$ws.Range("D27").Value = 799.207
$ws.Range("D28").Value = 875.065
$ws.Range("D29").Value = 961.752
$ws.Range("D30").Value = 1049.804
